$d = $word.ActiveDocument

# Split the single opening paragraph ("Github tutorial word doc") into four
# paragraphs:
#   1. "Github tutorial word doc"
#   2. (empty)
#   3. "Version 2"
#   4. (empty, holds the original _GoBack bookmark)
#
# Insert a temporary marker "X" on the blank line so the paragraph mark
# survives the Find/Replace as a normal line, then strip the marker back out
# to leave a truly empty <w:p/>.
$d.Content.Find.Execute("Github tutorial word doc", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Github tutorial word doc^pX^pVersion 2^p", 2)

$d.Content.Find.Execute("X", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)
